# Generate Report for Handback
#
# The handback transform for the "35b5f849-80b8-4a45-bca3-d724d281e62d"
# handoff failed (target file name returned by the transform,
# "briidrer.lsl", did not match the expected handoff-derived name) for
# both the zh-cn and de-de targets. Update the localization-status report:
#
#   - Overview sheet: roll-up Status for that file flips from
#     "Ready for handoff" to "Handback transform failed" for both the
#     zh-cn and de-de columns.
#   - zh-cn / de-de detail sheets: the same row's Status (column C) flips
#     the same way, and the Error Detail (column K) is populated with the
#     handback/handoff file-name-mismatch message for that language.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

# --- Overview roll-up (row 3 = 35b5f849-80b8-4a45-bca3-d724d281e62d.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusFailed
$overview.Range("C3").Value = $statusFailed

# --- Per-language detail sheets ---
$languages = @(
    @{ Sheet = "zh-cn"; Target = "35b5f849-80b8-4a45-bca3-d724d281e62d.bc843fc323f170d0cb0845f523c498cab111541f.zh-cn" },
    @{ Sheet = "de-de"; Target = "35b5f849-80b8-4a45-bca3-d724d281e62d.bc843fc323f170d0cb0845f523c498cab111541f.de-de" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $ws.Range("C3").Value = $statusFailed
    $ws.Range("K3").Value = "Handback file name: briidrer.lsl is different with handoff file name: $($lang.Target)."
}
